$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")

# The formula in D8 pulls the serialization path from the external
# FixedIncome.xla add-in (via the "Menu" external workbook link):
#   =[1]!qlSerializationPath(Trigger)
# Replace it with a plain literal path so the workbook no longer depends
# on that formula, then break/remove the external workbook link itself.
# (Leading "'" keeps Excel from re-interpreting the text and preserves the
# cell's existing text style/format rather than minting a new one.)
$ws.Range("D8").Value = "'C:\Users\erik\junk\"

# Remove the external reference to FixedIncome.xla's Menu.xla addin so the
# workbook no longer links out to it.
foreach ($link in @($wb.LinkSources())) {
    $wb.BreakLink($link, 1)
}
